$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old demo contents (days-of-week header + number grid) so the
# sheet starts clean before the new schedule data goes in.
$ws.Cells.Clear()

# New header row: "Time" / "Event" (columns C1:D1 keep a styled-but-empty
# placeholder cell, mirroring the old B1:D1 vertical-center styling)
$ws.Range("A1").Value = "Time"
$ws.Range("B1").Value = "Event"
$ws.Range("C1:D1").VerticalAlignment = -4108

# Schedule rows
$ws.Range("A2").Value = "Sat Dec 20 11:45:56 2025"
$ws.Range("B2").Value = "Basketball"

$ws.Range("A3").Value = "Sat Dec 20 11:46:05 2025"
$ws.Range("B3").Value = "Voleyball"

$ws.Range("A4").Value = "Sat Dec 20 11:46:11 2025"
$ws.Range("B4").Value = "Football"

$ws.Range("A5").Value = "Sat Dec 20 11:46:14 2025"
$ws.Range("B5").Value = "  "

$ws.Range("A6").Value = "Sat Dec 20 11:46:26 2025"
$ws.Range("B6").Value = "Note writing"

$ws.Range("A7").Value = "Sat Dec 20 11:46:33 2025"
$ws.Range("B7").Value = "CSC311"

$ws.Range("A8").Value = "March 1st, 2026"
$ws.Range("B8").Value = "CSC321"

$ws.Range("A9").Value = "Feb. 14th, 2026"
$ws.Range("B9").Value = "Life coaching"

# Column widths for the re-purposed Time/Event/.. columns
$ws.Columns.Item(1).ColumnWidth = 28.8
$ws.Columns.Item(2).ColumnWidth = 22.8
$ws.Columns.Item(3).ColumnWidth = 20.8
$ws.Columns.Item(4).ColumnWidth = 11.8

# Leave the selection on A3, as in the saved workbook
$ws.Range("A3").Select()

$wb.Save()
